$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 7.645999999999999
$ws.Range("B6").Value = 6.637
$ws.Range("B7").Value = 5.05
$ws.Range("C7").Value = -13.601
$ws.Range("C12").Value = -11.285
$ws.Range("E12").Value = 17.309
$ws.Range("D13").Value = -7.365
$ws.Range("D14").Value = -7.773999999999999
$ws.Range("C15").Value = -13.609
$ws.Range("B16").Value = 5.79
$ws.Range("D16").Value = -7.780000000000001
$ws.Range("D19").Value = -8.184999999999999
$ws.Range("B20").Value = 8.888999999999999
$ws.Range("C20").Value = -12.148
$ws.Range("C21").Value = -12.182
$ws.Range("C22").Value = -13.252
$ws.Range("D22").Value = -7.407999999999999
$ws.Range("E22").Value = 16.828
$ws.Range("C23").Value = -12.451
$ws.Range("B28").Value = 6.465000000000001
$ws.Range("B29").Value = 5.053
$ws.Range("C29").Value = -11.171
$ws.Range("E29").Value = 17.373
$ws.Range("B32").Value = 6.431
$ws.Range("C34").Value = -12.905
$ws.Range("E34").Value = 16.868
$ws.Range("D36").Value = -8.050000000000001
$ws.Range("B40").Value = 9.263999999999999
$ws.Range("C42").Value = -12.092
$ws.Range("C43").Value = -13.499
$ws.Range("E43").Value = 16.528
$ws.Range("C44").Value = -13.636
$ws.Range("C45").Value = -13.376
$ws.Range("B46").Value = 5.545
$ws.Range("C46").Value = -14.141
$ws.Range("D46").Value = -8.370000000000001
$ws.Range("E48").Value = 16.938
$ws.Range("C50").Value = -13.625
$ws.Range("D50").Value = -8.518000000000001
$ws.Range("B51").Value = 5.415000000000001
$ws.Range("C51").Value = -12.057
$ws.Range("B52").Value = 5.442
$ws.Range("B57").Value = 5.763000000000001
$ws.Range("B59").Value = 5.220999999999999
$ws.Range("E60").Value = 16.49
$ws.Range("B62").Value = 5.902
$ws.Range("B66").Value = 4.961
$ws.Range("C66").Value = -10.897
$ws.Range("C67").Value = -11.812
$ws.Range("E68").Value = 17.073
$ws.Range("E70").Value = 17.469
$ws.Range("B73").Value = 7.026999999999999
$ws.Range("E73").Value = 16.82
$ws.Range("B74").Value = 9.133000000000001
$ws.Range("C79").Value = -12.204
$ws.Range("C84").Value = -13.68
$ws.Range("E87").Value = 16.679
$ws.Range("B92").Value = 5.502
$ws.Range("C92").Value = -12.109
$ws.Range("E92").Value = 17.147
$ws.Range("D95").Value = -7.536
$ws.Range("C97").Value = -12.085
$ws.Range("D97").Value = -7.899999999999999
$ws.Range("B100").Value = 6.098000000000001
$ws.Range("E101").Value = 16.413

Write-Output "Updated 65 KNN-imputed cell values in Sheet1."
